$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): row5 F (想去人数) 4221 -> 4314, row7 F 55 -> 56
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 4314
$wsExhibit.Range("F7").Value = 56

# Sheet "全部类型" (all types): row9 F 4221 -> 4314, row11 F 55 -> 56
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 4314
$wsAll.Range("F11").Value = 56
